$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 ("Data"): prepend three new year rows (2024, 2023, 2022) ---
# Existing rows (header + years 2019..2006) all shift down by 3; their
# values are unchanged. Insert 3 blank rows right after the header row.
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).Insert()

# Column A holds years as text (shared strings), like the rest of the
# column, so force text storage (NumberFormat "@") then restore the
# default "Normal" style so no stray number format lingers on the cell.
$ws1.Cells.Item(2,1).NumberFormat = "@"
$ws1.Cells.Item(2,1).Value = "2024"
$ws1.Cells.Item(2,1).Style = "Normal"
$ws1.Cells.Item(2,2).Value = 0.8

$ws1.Cells.Item(3,1).NumberFormat = "@"
$ws1.Cells.Item(3,1).Value = "2023"
$ws1.Cells.Item(3,1).Style = "Normal"
$ws1.Cells.Item(3,2).Value = 1.2

$ws1.Cells.Item(4,1).NumberFormat = "@"
$ws1.Cells.Item(4,1).Value = "2022"
$ws1.Cells.Item(4,1).Style = "Normal"
$ws1.Cells.Item(4,2).Value = 0.62

# --- Sheet 2 ("Metadata"): update indicator text + add "actualizacion" row ---

# Row 1: key was an empty string; now a single space (merged with the
# existing " " string already used elsewhere on this sheet)
$ws2.Cells.Item(1,1).Value = " "
$ws2.Cells.Item(1,2).Value = " "

# nomindicador text rewritten
$ws2.Cells.Item(2,2).Value = "Porcentaje de personas que viven en viviendas sin agua potable"

# conindicador text rewritten
$ws2.Cells.Item(4,2).Value = "No acceso al agua potable"

# observaciones text rewritten (longer explanation replaces "Sin observaciones")
$observacionesText = @"
Desde marzo de 2020 hasta junio de 2021 se interrumpió el relevamiento presencial y se aplicó de manera telefónica un cuestionario restringido con el objetivo de continuar publicando los indicadores de ingresos y mercado de trabajo. En ese período la encuesta pasó a ser de paneles rotativos elegidos al azar a partir de los casos respondentes del año anterior. 
En julio de 2021 el INE retomó la realización de encuestas presenciales, pero introdujo un cambio metodológico, ya que la ECH pasa a ser una encuesta de panel rotativo con periodicidad mensual compuesta por seis paneles o grupos de rotación, cada uno de los cuales es una muestra representativa de la población. Con esta nueva metodología, cada hogar seleccionado participa durante seis meses de la ECH.
"@
$ws2.Cells.Item(8,2).Value = $observacionesText

# insert a new "actualizacion" row right after "observaciones" (row 8), before "cita"
$ws2.Rows.Item(9).Insert()
$ws2.Cells.Item(9,1).Value = "actualizacion"
$ws2.Cells.Item(9,2).Value = "Julio 2025"

# cita text rewritten (now row 10 after the insert), trailing newline included
$citaText = @"
UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE

"@
$ws2.Cells.Item(10,2).Value = $citaText
